$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- In-place value edits (rows 2-25), swapping which cells are "missing" ---
$ws.Range("D2").ClearContents()
$ws.Range("C6").Value = 15.1
$ws.Range("C8").ClearContents()
$ws.Range("C18").Value = 11.5
$ws.Range("C20").ClearContents()
$ws.Range("C23").Value = 12.2
$ws.Range("C25").ClearContents()

# --- Remove two records entirely: "RM 232" (row 26) and "SC 92" (originally row 28) ---
# After deleting row 26, the former row 28 ("SC 92") becomes row 27.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# --- Update remaining rows' missing/known values to match the new target layout ---
# Row 27 (SC 101): B becomes known
$ws.Range("B27").Value = -20.4

# Row 28 (SC 105): B becomes missing
$ws.Range("B28").ClearContents()

# Row 29 (SC 119): B becomes missing
$ws.Range("B29").ClearContents()

# Row 30 (SC 120): B, C, D become known
$ws.Range("B30").Value = -19.7
$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6

# Row 32 (SC 193): B becomes missing
$ws.Range("B32").ClearContents()
